# Auto-generated script to update market price-derived values across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1836.3438
$ws.Range("I15").Value = 1836.3438
$ws.Range("K15").Value = 5509.0314
$ws.Range("M15").Value = -5340.0314
$ws.Range("H32").Value = 2221.5293
$ws.Range("I32").Value = 1783.1666
$ws.Range("K32").Value = 1783.1666
$ws.Range("M32").Value = -1457.1666
$ws.Range("H70").Value = 2346.7
$ws.Range("J70").Value = 2645.875
$ws.Range("L70").Value = 7937.625
$ws.Range("N70").Value = -8477.625
$ws.Range("H73").Value = 2346.7
$ws.Range("J73").Value = 2645.875
$ws.Range("L73").Value = 7937.625
$ws.Range("N73").Value = -9809.625
$ws.Range("H82").Value = 8050.25
$ws.Range("I82").Value = 6200.4287
$ws.Range("K82").Value = 18601.2861
$ws.Range("M82").Value = -18195.2861
$ws.Range("H85").Value = 8050.25
$ws.Range("I85").Value = 6200.4287
$ws.Range("K85").Value = 18601.2861
$ws.Range("M85").Value = -17197.2861
$ws.Range("H92").Value = 67067.8
$ws.Range("I92").Value = 77078.30499999999
$ws.Range("J92").Value = 1999.5
$ws.Range("K92").Value = 77078.30499999999
$ws.Range("L92").Value = 1999.5
$ws.Range("M92").Value = -75830.30499999999
$ws.Range("N92").Value = -4495.5
$ws.Range("H132").Value = 3226.85
$ws.Range("I132").Value = 2463.4722
$ws.Range("J132").Value = 10097.25
$ws.Range("K132").Value = 7390.4166
$ws.Range("L132").Value = 30291.75
$ws.Range("M132").Value = -4860.4166
$ws.Range("N132").Value = -35351.75
$ws.Range("H137").Value = 4831.933
$ws.Range("I137").Value = 1250
$ws.Range("J137").Value = 8925.571
$ws.Range("K137").Value = 3750
$ws.Range("L137").Value = 26776.713
$ws.Range("M137").Value = -1200
$ws.Range("N137").Value = -31876.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30883.41
$ws.Range("I32").Value = 33278.03
$ws.Range("J32").Value = 14600
$ws.Range("K32").Value = 33278.03
$ws.Range("L32").Value = 14600
$ws.Range("M32").Value = -32991.03
$ws.Range("N32").Value = -15174
$ws.Range("H61").Value = 4383.1
$ws.Range("I61").Value = 4383.1
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4383.1
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4171.1
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 75321.78999999999
$ws.Range("I74").Value = 75321.78999999999
$ws.Range("K74").Value = 75321.78999999999
$ws.Range("M74").Value = -74447.78999999999
$ws.Range("H77").Value = 75321.78999999999
$ws.Range("I77").Value = 75321.78999999999
$ws.Range("K77").Value = 376608.95
$ws.Range("M77").Value = -372240.95
$ws.Range("H110").Value = 2567.85
$ws.Range("I110").Value = 2666.125
$ws.Range("J110").Value = 2174.75
$ws.Range("K110").Value = 2666.125
$ws.Range("L110").Value = 2174.75
$ws.Range("M110").Value = -621.125
$ws.Range("N110").Value = -6264.75
$ws.Range("H136").Value = 4383.1
$ws.Range("I136").Value = 4383.1
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13149.3
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10599.3
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2067.6365
$ws.Range("I86").Value = 1820.4615
$ws.Range("K86").Value = 1820.4615
$ws.Range("M86").Value = -697.4614999999999
$ws.Range("H89").Value = 2067.6365
$ws.Range("I89").Value = 1820.4615
$ws.Range("K89").Value = 9102.307499999999
$ws.Range("M89").Value = -3486.307499999999
$ws.Range("H94").Value = 9966.786
$ws.Range("I94").Value = 12223.5
$ws.Range("K94").Value = 12223.5
$ws.Range("M94").Value = -11772.5
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H134").Value = 2463.8333
$ws.Range("I134").Value = 2487.818
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 7463.454000000001
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -4928.454000000001
$ws.Range("N134").Value = -11670
$ws.Range("H141").Value = 80600
$ws.Range("J141").Value = 80600
$ws.Range("L141").Value = 80600
$ws.Range("N141").Value = -90960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 74872.5
$ws.Range("I58").Value = 80324.234
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 80324.234
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -80121.234
$ws.Range("N58").Value = -4406
$ws.Range("H94").Value = 3375.25
$ws.Range("J94").Value = 3342.1667
$ws.Range("L94").Value = 3342.1667
$ws.Range("N94").Value = -4244.1667
$ws.Range("H105").Value = 1311.4445
$ws.Range("J105").Value = 3794
$ws.Range("L105").Value = 3794
$ws.Range("N105").Value = -7288
$ws.Range("H132").Value = 1268
$ws.Range("I132").Value = 1162.6522
$ws.Range("J132").Value = 2883.3333
$ws.Range("K132").Value = 3487.9566
$ws.Range("L132").Value = 8649.999899999999
$ws.Range("M132").Value = -957.9566
$ws.Range("N132").Value = -13709.9999
$ws.Range("H136").Value = 74872.5
$ws.Range("I136").Value = 80324.234
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 240972.702
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -238422.702
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 381.75
$ws.Range("I68").Value = 479.25
$ws.Range("K68").Value = 1437.75
$ws.Range("M68").Value = -626.75
$ws.Range("H71").Value = 381.75
$ws.Range("I71").Value = 479.25
$ws.Range("K71").Value = 4313.25
$ws.Range("M71").Value = -257.25
$ws.Range("H107").Value = 716.5
$ws.Range("J107").Value = 327
$ws.Range("L107").Value = 981
$ws.Range("N107").Value = -4821

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H126").Value = 18331.666
$ws.Range("I126").Value = 14997.5
$ws.Range("K126").Value = 44992.5
$ws.Range("M126").Value = -42522.5
$ws.Range("H132").Value = 51979.05
$ws.Range("I132").Value = 73398.64
$ws.Range("K132").Value = 220195.92
$ws.Range("M132").Value = -217665.92

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4383.3335
$ws.Range("J68").Value = 6700
$ws.Range("L68").Value = 6700
$ws.Range("N68").Value = -8198
$ws.Range("H71").Value = 4383.3335
$ws.Range("J71").Value = 6700
$ws.Range("L71").Value = 33500
$ws.Range("N71").Value = -40988
$ws.Range("H82").Value = 3065.65
$ws.Range("J82").Value = 3147.9473
$ws.Range("L82").Value = 3147.9473
$ws.Range("N82").Value = -3869.9473
$ws.Range("H85").Value = 3065.65
$ws.Range("J85").Value = 3147.9473
$ws.Range("L85").Value = 3147.9473
$ws.Range("N85").Value = -5643.9473
$ws.Range("H93").Value = 1373.8182
$ws.Range("I93").Value = 781.6
$ws.Range("J93").Value = 2642.8572
$ws.Range("K93").Value = 781.6
$ws.Range("L93").Value = 2642.8572
$ws.Range("M93").Value = 466.4
$ws.Range("N93").Value = -5138.8572
$ws.Range("H133").Value = 67990.836
$ws.Range("J133").Value = 67990.836
$ws.Range("L133").Value = 67990.836
$ws.Range("N133").Value = -73050.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9241.5
$ws.Range("J81").Value = 12247.5
$ws.Range("L81").Value = 24495
$ws.Range("N81").Value = -26617
$ws.Range("H84").Value = 9241.5
$ws.Range("J84").Value = 12247.5
$ws.Range("L84").Value = 122475
$ws.Range("N84").Value = -133083
$ws.Range("H136").Value = 2053.6099
$ws.Range("I136").Value = 1727.909
$ws.Range("J136").Value = 3397.125
$ws.Range("K136").Value = 5183.727000000001
$ws.Range("L136").Value = 10191.375
$ws.Range("M136").Value = -2633.727000000001
$ws.Range("N136").Value = -15291.375
